$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.501.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.072.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.09"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.18%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +3.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0764"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.379.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.779"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.079.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.481.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +15.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0815"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.50%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.52%  "
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0625"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("E35").Value = "  +5.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.15%  "
$ws.Range("E37").Value = "  +5.16%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +23.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0952"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.462.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  +5.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "95.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.63%  "
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.71%  "
$ws.Range("E51").Value = "  +1.74%  "
